$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet: "query (44)" -> "query" ---
$ws.Name = "query"

# --- Insert 6 new data rows (188-193), copying formatting from the last existing data row (187) ---
$ws.Rows.Item(187).Copy()
$ws.Rows.Item(188).Insert()
$ws.Rows.Item(187).Copy()
$ws.Rows.Item(189).Insert()
$ws.Rows.Item(187).Copy()
$ws.Rows.Item(190).Insert()
$ws.Rows.Item(187).Copy()
$ws.Rows.Item(191).Insert()
$ws.Rows.Item(187).Copy()
$ws.Rows.Item(192).Insert()
$ws.Rows.Item(187).Copy()
$ws.Rows.Item(193).Insert()

# --- Populate the new rows with their values ---
# Row 188
$ws.Range("A188").Value = "Thiago"
$ws.Range("B188").Value = "Dúvida"
$ws.Range("C188").Value = "Cadastro VIP"
$ws.Range("D188").Value = 5459853000100
$ws.Range("E188").Value = "Elite Centro Automotivo De Serv Ltda"
$ws.Range("F188").Value = "Proprietaria pediu apoio para limpeza de base de vips antigos/deligado e apoio para inserção de 5 novos"
$ws.Range("G188").Value = "SIM"
$ws.Range("H188").Value = "Torre de Expansão"
$ws.Range("I188").Value = "Whatsapp"
$ws.Range("J188").Value = "Receptivo"
$ws.Range("K188").Value = 46010
$ws.Range("L188").Value = 46010
$ws.Range("M188").Value = "G.N. Urbano Ribeirao Preto"
$ws.Range("N188").Value = "Revendedor"
$ws.Range("O188").ClearContents()
$ws.Range("P188").ClearContents()
$ws.Range("Q188").Value = 0
$ws.Range("R188").ClearContents()
$ws.Range("S188").Value = "Item"
$ws.Range("T188").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(188).RowHeight = 29

# Row 189
$ws.Range("A189").Value = "Ricardo"
$ws.Range("B189").Value = "Dúvida"
$ws.Range("C189").Value = "Tratativas Financeiras/Repasses"
$ws.Range("D189").Value = 1332868000151
$ws.Range("E189").Value = "Brixner & Brixner Ltda"
$ws.Range("F189").Value = "O consultor Jairo entrou em contato para questionar o motivo pelo qual a revenda não está tendo acesso às NFs de MDR. Expliquei que já temos um chamado aberto sobre isso e que o time responsável já está tratando o caso."
$ws.Range("G189").Value = "NÃO"
$ws.Range("H189").Value = "Coord. De Controladoria"
$ws.Range("I189").Value = "Whatsapp"
$ws.Range("J189").Value = "Receptivo"
$ws.Range("K189").Value = 46010
$ws.Range("L189").Value = 46010
$ws.Range("M189").Value = "G.N. Urbano Londrina"
$ws.Range("N189").Value = "CN"
$ws.Range("O189").ClearContents()
$ws.Range("P189").ClearContents()
$ws.Range("Q189").Value = 0
$ws.Range("R189").ClearContents()
$ws.Range("S189").Value = "Item"
$ws.Range("T189").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(189).RowHeight = 43.5

# Row 190
$ws.Range("A190").Value = "Ricardo"
$ws.Range("B190").Value = "Aquisição"
$ws.Range("C190").Value = "KMV"
$ws.Range("D190").Value = 94094679000179
$ws.Range("E190").Value = "Arosi Comercio De Combustiveis Ltda"
$ws.Range("F190").Value = "Entrei em contato com o revendedor Pietro para passar as instruções iniciais sobre como operacionalizar o KMV, apresentar os conceitos e o Portal Parceiros KMV para controle financeiro."
$ws.Range("G190").Value = "SIM"
$ws.Range("H190").Value = "Torre de Expansão"
$ws.Range("I190").Value = "Telefone"
$ws.Range("J190").Value = "Ativo"
$ws.Range("K190").Value = 46010
$ws.Range("L190").Value = 46010
$ws.Range("M190").Value = "G.N. Urbano Santa Maria"
$ws.Range("N190").Value = "CN"
$ws.Range("O190").ClearContents()
$ws.Range("P190").ClearContents()
$ws.Range("Q190").Value = 0
$ws.Range("R190").ClearContents()
$ws.Range("S190").Value = "Item"
$ws.Range("T190").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(190).RowHeight = 29

# Row 191
$ws.Range("A191").Value = "Thiago"
$ws.Range("B191").Value = "Dúvida"
$ws.Range("C191").Value = "REDE"
$ws.Range("D191").Value = 12426856000102
$ws.Range("E191").Value = "Auto Posto Santa Luzita Ltda"
$ws.Range("F191").Value = "O Gerente Comercial solicitou apoio para intermediar com o adquirente a troca da maquininha, pois a atual estava obsoleta. Orientamos que o contato fosse feito diretamente com a central do próprio adquirente, já que não realizamos mais essa ponte. O gerente Willian efetuou a ligação e resolveu o caso. O novo equipamento já está instalado no posto."
$ws.Range("G191").Value = "SIM"
$ws.Range("H191").Value = "Revenda"
$ws.Range("I191").Value = "Whatsapp"
$ws.Range("J191").Value = "Receptivo"
$ws.Range("K191").Value = 46010
$ws.Range("L191").Value = 46010
$ws.Range("M191").Value = "G.N. Urbano Sp Leste"
$ws.Range("N191").Value = "Revendedor"
$ws.Range("O191").ClearContents()
$ws.Range("P191").ClearContents()
$ws.Range("Q191").Value = 0
$ws.Range("R191").ClearContents()
$ws.Range("S191").Value = "Item"
$ws.Range("T191").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(191).RowHeight = 58

# Row 192
$ws.Range("A192").Value = "Thiago"
$ws.Range("B192").Value = "Dúvida"
$ws.Range("C192").Value = "Dúvidas sobre conciliação"
$ws.Range("D192").Value = 50230537000116
$ws.Range("E192").Value = "Auto Posto Da Nova Jaboticabal Ltda"
$ws.Range("F192").Value = "Posto solicitou esclarecimentos sobre a taxa aplicada. Informei a taxa atual e como eles mesmos podem se autoatenderem no portal e verficarem por lá"
$ws.Range("G192").Value = "SIM"
$ws.Range("H192").Value = "Torre de Expansão"
$ws.Range("I192").Value = "Whatsapp"
$ws.Range("J192").Value = "Receptivo"
$ws.Range("K192").Value = 46010
$ws.Range("L192").Value = 46010
$ws.Range("M192").Value = "G.N. Urbano Ribeirao Preto"
$ws.Range("N192").Value = "Revendedor"
$ws.Range("O192").ClearContents()
$ws.Range("P192").ClearContents()
$ws.Range("Q192").Value = 0
$ws.Range("R192").ClearContents()
$ws.Range("S192").Value = "Item"
$ws.Range("T192").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(192).RowHeight = 29

# Row 193
$ws.Range("A193").Value = "Ricardo"
$ws.Range("B193").Value = "Aquisição"
$ws.Range("C193").Value = "KMV"
$ws.Range("D193").Value = 60289999000162
$ws.Range("E193").Value = "Goldani & Lima Com De Combustiveis Ltda"
$ws.Range("F193").Value = "Entrei em contato com o revendedor Alex, visto que ele realizou a adesão ao Conecta, porém não realizou a adesão ao KMV. Passei o entendimento para ele."
$ws.Range("G193").Value = "SIM"
$ws.Range("H193").Value = "Torre de Expansão"
$ws.Range("I193").Value = "Whatsapp"
$ws.Range("J193").Value = "Ativo"
$ws.Range("K193").Value = 46010
$ws.Range("L193").Value = 46010
$ws.Range("M193").Value = "G.N. Urbano Porto Alegre"
$ws.Range("N193").Value = "CN"
$ws.Range("O193").ClearContents()
$ws.Range("P193").ClearContents()
$ws.Range("Q193").Value = 0
$ws.Range("R193").ClearContents()
$ws.Range("S193").Value = "Item"
$ws.Range("T193").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(193).RowHeight = 29

# --- Update the hidden defined name (local part + range) ---
$definedName = $wb.Names.Item(1)
$definedName.Name = "query__46"
$definedName.RefersTo = "=query!`$A`$1:`$T`$193"

# --- Rename & resize the worksheet table (ListObject) ---
$lo = $ws.ListObjects.Item(1)
$lo.Name = "Tabela_query__46"
$lo.Resize($ws.Range("A1:T193"))

# --- Refresh the visible selection to match the new used range ---
$ws.Range("A1:T193").Select()

